# Add the two new "Follow-Up" header columns to the canned
# invalid-monitorees import/export template (Sara-Alert-Format-Invalid-Monitorees.xlsx).
#
# The sheet is a single header row (row 1) that lists every importable /
# exportable monitoree field. Two new fields were appended at the end:
#   DI1 -> "Follow-Up Reason"
#   DJ1 -> "Follow-Up Note"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing header row runs through column 112 (DH). Append the two new
# headers right after it.
$lastCol = 112
$followUpReasonCol = $lastCol + 1   # 113 -> DI
$followUpNoteCol   = $lastCol + 2   # 114 -> DJ

$ws.Cells.Item(1, $followUpReasonCol).Value = "Follow-Up Reason"
$ws.Cells.Item(1, $followUpNoteCol).Value = "Follow-Up Note"

# Match the bespoke column widths the workbook author set for these two
# new header cells (values are expressed in Excel's character-width units,
# same scale used by the rest of the sheet's "bestFit" columns).
$ws.Columns.Item($followUpReasonCol).ColumnWidth = 14
$ws.Columns.Item($followUpNoteCol).ColumnWidth = 12.17
